# Updates cryptos list values (price & volume%) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.002.75"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.343.53"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.17"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.45"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.342.30"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.736.70"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.45"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.066.51"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.339.68"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.54"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.42"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.80"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.71"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.46"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0724"
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.53"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.918"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.20"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.58"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "280.53"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.07"
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0930"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("E47").Value = "  -1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.90"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("E50").Value = "  -2.03%  "
